$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 534.5
$ws.Range("J17").Value = 545.03705
$ws.Range("L17").Value = 1635.11115
$ws.Range("N17").Value = -1971.11115
$ws.Range("H32").Value = 4796.5
$ws.Range("J32").Value = 2189.4443
$ws.Range("L32").Value = 2189.4443
$ws.Range("N32").Value = -2841.4443
$ws.Range("H33").Value = 469.14285
$ws.Range("I33").Value = 222.25
$ws.Range("K33").Value = 222.25
$ws.Range("M33").Value = 6.75
$ws.Range("H40").Value = 4489.125
$ws.Range("I40").Value = 4615.6665
$ws.Range("K40").Value = 4615.6665
$ws.Range("M40").Value = -4440.6665
$ws.Range("H69").Value = 166784.47
$ws.Range("J69").Value = 191289.77
$ws.Range("L69").Value = 573869.3099999999
$ws.Range("N69").Value = -575617.3099999999
$ws.Range("H70").Value = 3314.6785
$ws.Range("I70").Value = 1722.6364
$ws.Range("K70").Value = 5167.9092
$ws.Range("M70").Value = -4897.9092
$ws.Range("H72").Value = 166784.47
$ws.Range("J72").Value = 191289.77
$ws.Range("L72").Value = 1721607.93
$ws.Range("N72").Value = -1730343.93
$ws.Range("H73").Value = 3314.6785
$ws.Range("I73").Value = 1722.6364
$ws.Range("K73").Value = 5167.9092
$ws.Range("M73").Value = -4231.9092
$ws.Range("H96").Value = 4464856.5
$ws.Range("I96").Value = 7936778.5
$ws.Range("K96").Value = 23810335.5
$ws.Range("M96").Value = -23808962.5
$ws.Range("H101").Value = 238
$ws.Range("I101").Value = 253.2
$ws.Range("K101").Value = 759.5999999999999
$ws.Range("M101").Value = 862.4000000000001
$ws.Range("H116").Value = 18920.84
$ws.Range("I116").Value = 21024.316
$ws.Range("J116").Value = 12259.833
$ws.Range("K116").Value = 21024.316
$ws.Range("L116").Value = 12259.833
$ws.Range("M116").Value = -17582.316
$ws.Range("N116").Value = -19143.833
$ws.Range("H138").Value = 24264.195
$ws.Range("J138").Value = 31861.205
$ws.Range("L138").Value = 95583.61500000001
$ws.Range("N138").Value = -105863.615

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 790.3570999999999
$ws.Range("I2").Value = 776
$ws.Range("K2").Value = 776
$ws.Range("M2").Value = -663
$ws.Range("H32").Value = 19134.352
$ws.Range("I32").Value = 20495.793
$ws.Range("J32").Value = 1095.25
$ws.Range("K32").Value = 20495.793
$ws.Range("L32").Value = 1095.25
$ws.Range("M32").Value = -20208.793
$ws.Range("N32").Value = -1669.25
$ws.Range("H45").Value = 3275.9375
$ws.Range("I45").Value = 1901.7
$ws.Range("K45").Value = 1901.7
$ws.Range("M45").Value = -1524.7
$ws.Range("H61").Value = 7559.6665
$ws.Range("I61").Value = 1026.1
$ws.Range("K61").Value = 1026.1
$ws.Range("M61").Value = -814.0999999999999
$ws.Range("H110").Value = 673.25
$ws.Range("I110").Value = 673.25
$ws.Range("K110").Value = 673.25
$ws.Range("M110").Value = 1371.75
$ws.Range("H116").Value = 790.3570999999999
$ws.Range("I116").Value = 776
$ws.Range("K116").Value = 776
$ws.Range("M116").Value = 1518
$ws.Range("H132").Value = 1550.75
$ws.Range("I132").Value = 1185.9
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 3557.7
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -1027.7
$ws.Range("N132").Value = -15185
$ws.Range("H136").Value = 7559.6665
$ws.Range("I136").Value = 1026.1
$ws.Range("K136").Value = 3078.3
$ws.Range("M136").Value = -528.2999999999997

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 790.3570999999999
$ws.Range("I3").Value = 776
$ws.Range("K3").Value = 776
$ws.Range("M3").Value = -662
$ws.Range("H86").Value = 1525.0555
$ws.Range("I86").Value = 1495.4
$ws.Range("K86").Value = 1495.4
$ws.Range("M86").Value = -372.4000000000001
$ws.Range("H89").Value = 1525.0555
$ws.Range("I89").Value = 1495.4
$ws.Range("K89").Value = 7477
$ws.Range("M89").Value = -1861

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2372
$ws.Range("I16").Value = 1518.4
$ws.Range("K16").Value = 1518.4
$ws.Range("M16").Value = -1231.4
$ws.Range("H31").Value = 5558008.5
$ws.Range("I31").Value = 7693285.5
$ws.Range("K31").Value = 7693285.5
$ws.Range("M31").Value = -7692990.5
$ws.Range("H34").Value = 5558008.5
$ws.Range("I34").Value = 7693285.5
$ws.Range("K34").Value = 7693285.5
$ws.Range("M34").Value = -7693083.5
$ws.Range("H99").Value = 6047.75
$ws.Range("I99").Value = 4268.923
$ws.Range("K99").Value = 4268.923
$ws.Range("M99").Value = -2770.923
$ws.Range("H107").Value = 592
$ws.Range("I107").Value = 475.04544
$ws.Range("K107").Value = 475.04544
$ws.Range("M107").Value = 1444.95456
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H113").Value = 2372
$ws.Range("I113").Value = 1518.4
$ws.Range("K113").Value = 1518.4
$ws.Range("M113").Value = 651.5999999999999
$ws.Range("H126").Value = 6047.75
$ws.Range("I126").Value = 4268.923
$ws.Range("K126").Value = 12806.769
$ws.Range("M126").Value = -10336.769
$ws.Range("H132").Value = 51634.1
$ws.Range("I132").Value = 84400.336
$ws.Range("J132").Value = 2484.75
$ws.Range("K132").Value = 253201.008
$ws.Range("L132").Value = 7454.25
$ws.Range("M132").Value = -250671.008
$ws.Range("N132").Value = -12514.25
$ws.Range("H134").Value = 2582.5386
$ws.Range("I134").Value = 2494.9
$ws.Range("K134").Value = 7484.700000000001
$ws.Range("M134").Value = -4949.700000000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 570.6667
$ws.Range("I134").Value = 570.6667
$ws.Range("K134").Value = 1712.0001
$ws.Range("M134").Value = 3357.9999
$ws.Range("H140").Value = 3401.0833
$ws.Range("I140").Value = 3401.0833
$ws.Range("K140").Value = 10203.2499
$ws.Range("M140").Value = -5023.249899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 552.875
$ws.Range("I2").Value = 763.0769
$ws.Range("J2").Value = 304.45456
$ws.Range("K2").Value = 763.0769
$ws.Range("L2").Value = 304.45456
$ws.Range("M2").Value = -650.0769
$ws.Range("N2").Value = -530.45456
$ws.Range("H55").Value = 5445
$ws.Range("J55").Value = 8440
$ws.Range("L55").Value = 8440
$ws.Range("N55").Value = -9094
$ws.Range("H126").Value = 2499.75
$ws.Range("I126").Value = 2499.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7499.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5029.25
$ws.Range("N126").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4610.6665
$ws.Range("I7").Value = 3916.3333
$ws.Range("J7").Value = 5999.3335
$ws.Range("K7").Value = 3916.3333
$ws.Range("L7").Value = 5999.3335
$ws.Range("M7").Value = -3804.3333
$ws.Range("N7").Value = -6223.3335
$ws.Range("H122").Value = 3040
$ws.Range("I122").Value = 3333.3333
$ws.Range("K122").Value = 9999.999899999999
$ws.Range("M122").Value = -7549.999899999999
$ws.Range("H126").Value = 4610.6665
$ws.Range("I126").Value = 3916.3333
$ws.Range("J126").Value = 5999.3335
$ws.Range("K126").Value = 11748.9999
$ws.Range("L126").Value = 17998.0005
$ws.Range("M126").Value = -9278.999899999999
$ws.Range("N126").Value = -22938.0005
$ws.Range("H132").Value = 3684.7646
$ws.Range("I132").Value = 3543.5454
$ws.Range("J132").Value = 3943.6667
$ws.Range("K132").Value = 10630.6362
$ws.Range("L132").Value = 11831.0001
$ws.Range("M132").Value = -8100.636200000001
$ws.Range("N132").Value = -16891.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1360.6666
$ws.Range("I107").Value = 1624.8334
$ws.Range("J107").Value = 1096.5
$ws.Range("K107").Value = 4874.5002
$ws.Range("L107").Value = 3289.5
$ws.Range("M107").Value = -2954.5002
$ws.Range("N107").Value = -7129.5
$ws.Range("H132").Value = 60225.25
$ws.Range("I132").Value = 78523.11
$ws.Range("K132").Value = 235569.33
$ws.Range("M132").Value = -233039.33
$ws.Range("H136").Value = 22534
$ws.Range("I136").Value = 29743.2
$ws.Range("J136").Value = 4511
$ws.Range("K136").Value = 89229.60000000001
$ws.Range("L136").Value = 13533
$ws.Range("M136").Value = -86679.60000000001
$ws.Range("N136").Value = -18633

Write-Host "Applied Midgardsormr Profits updates"